# Apply the changes described by the diff:
# - Summary!B4: 100 -> 50
# - Symbols!A2: AAPL -> SPY
# - Symbols!B2: 100 -> 50
# - Symbols!E2: "Options on Apple" -> "Options on S&P 500 ETF"
# - Strategies!A2: AAPL -> SPY
# - Strategies!B2: "2025-02-21" -> "2025-01-10" (must stay plain text, not become a date serial)
# - Strategies!C2: "Long Call" -> "Short Put"
# - Strategies!D2: 100 -> 50
# - Strategies!E2: 1 -> 0.08333333333333333
# - Strategies!F2: 100 -> 600
# - Strategies!G2: "Options on Apple" -> "Options on S&P 500 ETF"

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B4").Value = 50

# --- Symbols sheet ---
$symbols = $wb.Worksheets.Item("Symbols")
$symbols.Range("A2").Value = "SPY"
$symbols.Range("B2").Value = 50
$symbols.Range("E2").Value = "Options on S&P 500 ETF"

# --- Strategies sheet ---
$strategies = $wb.Worksheets.Item("Strategies")
$strategies.Range("A2").Value = "SPY"

# B2 holds a date-looking string ("2025-01-10") that must remain plain text
# (it was stored as inlineStr/text before the edit, not a real date). A bare
# assignment gets auto-converted to a date serial by Excel's smart entry, so
# prefix with an apostrophe to force text, then reset the cell style back to
# Normal so the quote-prefix flag doesn't leave a stray style behind.
$strategies.Range("B2").Value = "'2025-01-10"
$strategies.Range("B2").Style = "Normal"

$strategies.Range("C2").Value = "Short Put"
$strategies.Range("D2").Value = 50
$strategies.Range("E2").Value = 0.08333333333333333
$strategies.Range("F2").Value = 600
$strategies.Range("G2").Value = "Options on S&P 500 ETF"
